$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5781.5557
$ws.Range("I19").Value = 887.5
$ws.Range("J19").Value = 9696.799999999999
$ws.Range("K19").Value = 887.5
$ws.Range("L19").Value = 9696.799999999999
$ws.Range("M19").Value = -712.5
$ws.Range("N19").Value = -10046.8
$ws.Range("H31").Value = 6995
$ws.Range("I31").Value = 6995
$ws.Range("K31").Value = 20985
$ws.Range("M31").Value = -20755
$ws.Range("H40").Value = 2878.5881
$ws.Range("I40").Value = 2811.75
$ws.Range("J40").Value = 3039
$ws.Range("K40").Value = 2811.75
$ws.Range("L40").Value = 3039
$ws.Range("M40").Value = -2636.75
$ws.Range("N40").Value = -3389
$ws.Range("H125").Value = 75658
$ws.Range("J125").Value = 800
$ws.Range("L125").Value = 7200
$ws.Range("N125").Value = -12120
$ws.Range("H137").Value = 13309.5
$ws.Range("J137").Value = 5727.143
$ws.Range("L137").Value = 17181.429
$ws.Range("N137").Value = -22281.429
$ws.Range("H138").Value = 20796.877
$ws.Range("I138").Value = 2328.7334
$ws.Range("J138").Value = 41317.035
$ws.Range("K138").Value = 6986.2002
$ws.Range("L138").Value = 123951.105
$ws.Range("M138").Value = -1846.2002
$ws.Range("N138").Value = -134231.105

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29559.406
$ws.Range("I32").Value = 22163.8
$ws.Range("K32").Value = 22163.8
$ws.Range("M32").Value = -21876.8
$ws.Range("H45").Value = 4666.3335
$ws.Range("I45").Value = 3027.4443
$ws.Range("J45").Value = 6305.222
$ws.Range("K45").Value = 3027.4443
$ws.Range("L45").Value = 6305.222
$ws.Range("M45").Value = -2650.4443
$ws.Range("N45").Value = -7059.222
$ws.Range("H61").Value = 4120.75
$ws.Range("I61").Value = 1233.75
$ws.Range("K61").Value = 1233.75
$ws.Range("M61").Value = -1021.75
$ws.Range("H88").Value = 12250
$ws.Range("I88").Value = 3000
$ws.Range("K88").Value = 3000
$ws.Range("M88").Value = -2594
$ws.Range("H91").Value = 12250
$ws.Range("I91").Value = 3000
$ws.Range("K91").Value = 3000
$ws.Range("M91").Value = -1596
$ws.Range("H132").Value = 1423.7
$ws.Range("I132").Value = 1162.7632
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 3488.2896
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -958.2896000000001
$ws.Range("N132").Value = -11810
$ws.Range("H136").Value = 4120.75
$ws.Range("I136").Value = 1233.75
$ws.Range("K136").Value = 3701.25
$ws.Range("M136").Value = -1151.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 26442.818
$ws.Range("I75").Value = 16174.6
$ws.Range("J75").Value = 34999.668
$ws.Range("K75").Value = 16174.6
$ws.Range("L75").Value = 34999.668
$ws.Range("M75").Value = -15238.6
$ws.Range("N75").Value = -36871.668
$ws.Range("H78").Value = 26442.818
$ws.Range("I78").Value = 16174.6
$ws.Range("J78").Value = 34999.668
$ws.Range("K78").Value = 48523.8
$ws.Range("L78").Value = 104999.004
$ws.Range("M78").Value = -43843.8
$ws.Range("N78").Value = -114359.004

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8021.5
$ws.Range("I62").Value = 7895.1763
$ws.Range("J62").Value = 8328.286
$ws.Range("K62").Value = 7895.1763
$ws.Range("L62").Value = 8328.286
$ws.Range("M62").Value = -7271.1763
$ws.Range("N62").Value = -9576.286
$ws.Range("H65").Value = 8021.5
$ws.Range("I65").Value = 7895.1763
$ws.Range("J65").Value = 8328.286
$ws.Range("K65").Value = 39475.8815
$ws.Range("L65").Value = 41641.43
$ws.Range("M65").Value = -36355.8815
$ws.Range("N65").Value = -47881.43

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 67.09090999999999
$ws.Range("I33").Value = 75.375
$ws.Range("J33").Value = 45
$ws.Range("K33").Value = 452.25
$ws.Range("L33").Value = 270
$ws.Range("M33").Value = -169.25
$ws.Range("N33").Value = -836
$ws.Range("H34").Value = 5682
$ws.Range("J34").Value = 8443.286
$ws.Range("L34").Value = 25329.858
$ws.Range("N34").Value = -25497.858
$ws.Range("H39").Value = 11252
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 13336
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 40008
$ws.Range("M39").Value = -14706
$ws.Range("N39").Value = -40596
$ws.Range("H51").Value = 3875.074
$ws.Range("H55").Value = 840
$ws.Range("I55").Value = 840
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 2520
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2343
$ws.Range("N55").ClearContents()
$ws.Range("H63").Value = 6609.2173
$ws.Range("I63").Value = 17337.334
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 52012.00199999999
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -51263.00199999999
$ws.Range("N63").Value = -16498
$ws.Range("H66").Value = 6609.2173
$ws.Range("I66").Value = 17337.334
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 156036.006
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -152292.006
$ws.Range("N66").Value = -52488
$ws.Range("H75").Value = 4597.727
$ws.Range("J75").Value = 5000
$ws.Range("L75").Value = 15000
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 4597.727
$ws.Range("J78").Value = 5000
$ws.Range("L78").Value = 45000
$ws.Range("N78").Value = -54984
$ws.Range("H113").Value = 733
$ws.Range("J113").Value = 899.5
$ws.Range("L113").Value = 2698.5
$ws.Range("N113").Value = -7038.5
$ws.Range("H131").Value = 1855.6296
$ws.Range("I131").Value = 1341.4667
$ws.Range("J131").Value = 2498.3333
$ws.Range("K131").Value = 4024.4001
$ws.Range("L131").Value = 7494.999899999999
$ws.Range("M131").Value = 1015.5999
$ws.Range("N131").Value = -17574.9999
$ws.Range("H139").Value = 4306.091
$ws.Range("J139").Value = 5000
$ws.Range("L139").Value = 15000
$ws.Range("N139").Value = -25280

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13510.111
$ws.Range("I70").Value = 13691
$ws.Range("K70").Value = 13691
$ws.Range("M70").Value = -13421
$ws.Range("H73").Value = 13510.111
$ws.Range("I73").Value = 13691
$ws.Range("K73").Value = 13691
$ws.Range("M73").Value = -12755
$ws.Range("H80").Value = 7989.476
$ws.Range("I80").Value = 4884.5557
$ws.Range("K80").Value = 4884.5557
$ws.Range("M80").Value = -3886.5557
$ws.Range("H83").Value = 7989.476
$ws.Range("I83").Value = 4884.5557
$ws.Range("K83").Value = 24422.7785
$ws.Range("M83").Value = -19430.7785
$ws.Range("H107").Value = 257.18182
$ws.Range("J107").Value = 319.8
$ws.Range("L107").Value = 319.8
$ws.Range("N107").Value = -4159.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 906.3333
$ws.Range("I61").Value = 780.6667
$ws.Range("J61").Value = 1283.3334
$ws.Range("K61").Value = 780.6667
$ws.Range("L61").Value = 1283.3334
$ws.Range("M61").Value = -578.6667
$ws.Range("N61").Value = -1687.3334
$ws.Range("H68").Value = 3258.7646
$ws.Range("I68").Value = 2741.4167
$ws.Range("J68").Value = 4500.4
$ws.Range("K68").Value = 2741.4167
$ws.Range("L68").Value = 4500.4
$ws.Range("M68").Value = -1992.4167
$ws.Range("N68").Value = -5998.4
$ws.Range("H71").Value = 3258.7646
$ws.Range("I71").Value = 2741.4167
$ws.Range("J71").Value = 4500.4
$ws.Range("K71").Value = 13707.0835
$ws.Range("L71").Value = 22502
$ws.Range("M71").Value = -9963.083500000001
$ws.Range("N71").Value = -29990
$ws.Range("H113").Value = 906.3333
$ws.Range("I113").Value = 780.6667
$ws.Range("J113").Value = 1283.3334
$ws.Range("K113").Value = 780.6667
$ws.Range("L113").Value = 1283.3334
$ws.Range("M113").Value = 1389.3333
$ws.Range("N113").Value = -5623.3334
$ws.Range("H122").Value = 4118.6855
$ws.Range("I122").Value = 3166.087
$ws.Range("K122").Value = 9498.261
$ws.Range("M122").Value = -7048.261
$ws.Range("H136").Value = 3016.6135
$ws.Range("I136").Value = 3033.2563
$ws.Range("J136").Value = 2886.8
$ws.Range("K136").Value = 9099.768899999999
$ws.Range("L136").Value = 8660.400000000001
$ws.Range("M136").Value = -6549.768899999999
$ws.Range("N136").Value = -13760.4

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 71000
$ws.Range("I51").Value = 51333.332
$ws.Range("K51").Value = 51333.332
$ws.Range("M51").Value = -50823.332
$ws.Range("H57").Value = 79999.5
$ws.Range("J57").Value = 79999.5
$ws.Range("L57").Value = 79999.5
$ws.Range("N57").Value = -81507.5
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 6430183
$ws.Range("I132").Value = 9643821
$ws.Range("K132").Value = 28931463
$ws.Range("M132").Value = -28928933
